$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts all existing columns
# (A:G -> B:H) right by one, preserving their values and styles.
$ws.Columns("A").Insert()

# Populate the new "Status" column (header + merge-result per row).
$ws.Range("A1").Value = "Status"
$ws.Range("A2").Value = "matched"
$ws.Range("A3").Value = "not matched"
$ws.Range("A4").Value = "not matched"
$ws.Range("A5").Value = "missing"
$ws.Range("A6").Value = "missing"

# Match the header's bold/border/alignment style used by the rest of row 1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
